$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 553.9339164698035
$ws.Range("D2").Value = 136.0796753216453
$ws.Range("G2").Value = 513
$ws.Range("C3").Value = 37.54371199562761
$ws.Range("D3").Value = 6.493415194419221
$ws.Range("F3").Value = 32.81
$ws.Range("G3").Value = 37.81
$ws.Range("H3").Value = 41.51
$ws.Range("C4").Value = 1.98216252087201
$ws.Range("D4").Value = 2.542316588930305
$ws.Range("G4").Value = 1.27
$ws.Range("H4").Value = 2.41
$ws.Range("C5").Value = 323.3209024386936
$ws.Range("D5").Value = 10.80527007418991
$ws.Range("F5").Value = 316.54
$ws.Range("G5").Value = 324.78
$ws.Range("H5").Value = 331.49
$ws.Range("C6").Value = 21.20733341915513
$ws.Range("D6").Value = 2.579251734900613
$ws.Range("F6").Value = 19.73
$ws.Range("G6").Value = 21.17
$ws.Range("H6").Value = 22.53
$ws.Range("C7").Value = -76.90280080160562
$ws.Range("D7").Value = 22.89926493007523
$ws.Range("C8").Value = 7.419475878633727
$ws.Range("D8").Value = 7.098372356560299
$ws.Range("C9").Value = 9.322680892004572
$ws.Range("D9").Value = 1.685704350771354
$ws.Range("C10").Value = 867.8301139623015
$ws.Range("D10").Value = 0.461504525285086
$ws.Range("C11").Value = 0.5559225975394744
$ws.Range("D11").Value = 0.5890156371207561
$ws.Range("C12").Value = 22.74623852133575
$ws.Range("D12").Value = 12.29406113203675
$ws.Range("C13").Value = 0.6740255086446632
$ws.Range("D13").Value = 0.7506961663209128
$ws.Range("C14").Value = 1.827532081539733
$ws.Range("D14").Value = 1.664426921769673
$ws.Range("C15").Value = 94.1628008016058
$ws.Range("D15").Value = 22.89926493007523
$ws.Range("C16").Value = -86.02480107908279
$ws.Range("D16").Value = 20.47195713111128
$ws.Range("F16").Value = -102.7376019773414
$ws.Range("G16").Value = -84.23249407632485
$ws.Range("H16").Value = -70.46683163887967
$ws.Range("C17").Value = -78.60532520044909
$ws.Range("D17").Value = 25.44454038663014
$ws.Range("F17").Value = -93.87736039420676
$ws.Range("G17").Value = -73.57382219273629
$ws.Range("H17").Value = -60.41392685158225
